$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price/volume figures pulled by the symbol-list refresh action.
# Column D values are stored as text in this sheet, so each cell is forced
# to a text number format before the assignment (otherwise a numeric-looking
# string like "243.50" would be auto-converted to the number 243.5) and the
# style is reset back to Normal afterwards so no stray formatting is left behind.
$priceUpdates = @{
    "D2" = "243.50"
    "D3" = "23.08"
    "D4" = "5.402"
    "D6" = "3.457"
    "D7" = "6.555"
    "D8" = "0.8120"
    "D10" = "0.1417"
    "D11" = "0.07396"
    "D12" = "0.03275"
    "D14" = "0.09352"
    "D15" = "3.855"
    "D16" = "0.001569"
    "D17" = "0.04681"
    "D18" = "0.0005943"
    "D19" = "0.005946"
    "D20" = "0.004953"
    "D21" = "0.0009865"
    "D22" = "0.00008602"
    "D23" = "3.610"
    "D40" = "0.03966"
    "D42" = "0.1076"
    "D44" = "0.008624"
    "D45" = "0.00005171"
    "D47" = "0.7203"
    "D48" = "0.002265"
    "D49" = "0.00002100"
    "D50" = "0.0002000"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Column E labels (plain text, no numeric coercion risk).
$labelUpdates = @{
    "E18" = "17OneONEWorstin24h"
    "E41" = "40KickTokenKICK"
    "E44" = "43LocalTradersLCTBestin24h"
    "E47" = "46CoinbaseStockTokenCOIN"
}

foreach ($addr in $labelUpdates.Keys) {
    $ws.Range($addr).Value = $labelUpdates[$addr]
}

